$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$win = $excel.ActiveWindow
try {
  $win.SmallScroll(5,0,0,0)
  Write-Host "SmallScroll OK"
} catch {
  Write-Host ("SmallScroll failed: " + $_.Exception.Message)
}
